$d = $word.ActiveDocument

# 1. Update the digit validation error message text.
[void]$d.Content.Find.Execute(
    "Password must contain at least one digit",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Password must contain at least one numeric digit", 2)

# 2. The blank paragraph that used to follow that line now holds a
#    run of four spaces instead of being fully empty.
$blank = $d.Paragraphs.Item(9)
$blank.Range.Text = "    "

# 3. Collapse the double blank paragraph before "if __name__ ==" (at
#    the end of the last print(valid, msg) block) down to a single
#    blank paragraph.
$extra = $d.Paragraphs.Item(30)
$extra.Range.Delete()
